# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the existing columns (B1:G1) and writing the
# numeric flag for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the
# new header cell (H1) so it picks up the same bold/border/alignment
# style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for row 2 — plain numeric value, no special style,
# matching the other numeric cells in the row.
$ws.Range("H2").Value = 1
